$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.232.66"
$ws.Range("E2").Value = "  +2.17%  "
$ws.Range("D3").Value = "2.362.62"
$ws.Range("E3").Value = "  +0.44%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.678"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.67%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.48"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +10.64%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.552"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +20.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.63"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +16.61%  "
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").Value = "2.711.75"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "16.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +10.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.909"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.27%  "
$ws.Range("D17").Value = "2.366.39"
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("D18").Value = "44.252.66"
$ws.Range("E18").Value = "  +2.25%  "
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.38%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "255.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.58%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.36%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "174.15"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.14%  "
$ws.Range("E30").Value = "  +3.25%  "
$ws.Range("E31").Value = "  +3.86%  "
$ws.Range("E32").Value = "  +5.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0747"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.61%  "
$ws.Range("E34").Value = "  +4.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.60%  "
$ws.Range("E37").Value = "  -1.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("E39").Value = "  +6.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.29"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.56%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  -1.13%  "
$ws.Range("E43").Value = "  +3.64%  "
$ws.Range("E44").Value = "  +4.54%  "
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.186"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "99.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.46"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.01%  "
$ws.Range("E49").Value = "  +5.13%  "
$ws.Range("D50").Value = "1.447.23"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.000205"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.05%  "
